# Insert two new rows at position 512 (pushing the existing rows 512-613
# down to 514-615) and populate them with the new week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("512:513").Insert()

# Row 512 - Pepino dulce, Cultivar IV Región, Primera
$ws.Range("A512").Value2 = 8
$ws.Range("B512").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C512").Value2 = "Coquimbo"
$ws.Range("D512").Value2 = 45015
$ws.Range("E512").Value2 = 4
$ws.Range("F512").Value2 = 100112043
$ws.Range("G512").Value2 = "Pepino dulce"
$ws.Range("H512").Value2 = "Cultivar IV Región"
$ws.Range("I512").Value2 = "Primera"
$ws.Range("J512").Value2 = 600
$ws.Range("K512").Value2 = 11000
$ws.Range("L512").Value2 = 12000
$ws.Range("M512").Value2 = 11500
$ws.Range("N512").Value2 = "$/bandeja 18 kilos"
$ws.Range("O512").Value2 = "Provincia de Limarí"
$ws.Range("P512").Value2 = 639
$ws.Range("Q512").Value2 = 18
$ws.Range("R512").Value2 = "Hortaliza"

# Row 513 - Pepino dulce, Cultivar IV Región, Segunda
$ws.Range("A513").Value2 = 8
$ws.Range("B513").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C513").Value2 = "Coquimbo"
$ws.Range("D513").Value2 = 45015
$ws.Range("E513").Value2 = 4
$ws.Range("F513").Value2 = 100112043
$ws.Range("G513").Value2 = "Pepino dulce"
$ws.Range("H513").Value2 = "Cultivar IV Región"
$ws.Range("I513").Value2 = "Segunda"
$ws.Range("J513").Value2 = 340
$ws.Range("K513").Value2 = 9000
$ws.Range("L513").Value2 = 10000
$ws.Range("M513").Value2 = 9500
$ws.Range("N513").Value2 = "$/bandeja 18 kilos"
$ws.Range("O513").Value2 = "Provincia de Limarí"
$ws.Range("P513").Value2 = 528
$ws.Range("Q513").Value2 = 18
$ws.Range("R513").Value2 = "Hortaliza"
